$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.071.26'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.837.12'
$ws.Range('E3').Value = '  +0.55%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.ClearFormats()
$ws.Range('E4').Value = '  +0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '242.99'
$c.ClearFormats()
$ws.Range('E5').Value = '  +0.74%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '0.6286'
$c.ClearFormats()
$ws.Range('E6').Value = '  -1.31%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range('E7').Value = '  +0.00%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.07589'
$c.ClearFormats()
$ws.Range('E8').Value = '  +3.41%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.2932'
$c.ClearFormats()
$ws.Range('E9').Value = '  -0.13%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '22.61'
$c.ClearFormats()
$ws.Range('E10').Value = '  -0.87%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.07749'
$c.ClearFormats()
$ws.Range('D12').Value = '1.841.67'
$ws.Range('E12').Value = '  +0.84%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '4.972'
$c.ClearFormats()
$ws.Range('E13').Value = '  -0.27%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.6663'
$c.ClearFormats()
$ws.Range('E14').Value = '  +0.37%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.00001003'
$c.ClearFormats()
$ws.Range('E15').Value = '  +15.60%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '83.18'
$c.ClearFormats()
$ws.Range('E16').Value = '  +1.43%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '6.066'
$c.ClearFormats()
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '29.090.74'
$ws.Range('E18').Value = '  +0.72%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '227.08'
$c.ClearFormats()
$ws.Range('E19').Value = '  +1.24%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '12.43'
$c.ClearFormats()
$ws.Range('E20').Value = '  +0.25%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.ClearFormats()
$ws.Range('E21').Value = '  +0.03%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '7.224'
$c.ClearFormats()
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.004'
$c.ClearFormats()
$ws.Range('E23').Value = '  +0.38%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '159.54'
$c.ClearFormats()
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('B25').Value = 'Stellar'
$ws.Range('C25').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.1387'
$c.ClearFormats()
$ws.Range('E25').Value = '  +1.24%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '8.509'
$c.ClearFormats()
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('E27').Value = '  +0.43%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.496'
$c.ClearFormats()
$ws.Range('E28').Value = '  -0.60%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '4.101'
$c.ClearFormats()
$ws.Range('E29').Value = '  +0.26%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.014'
$c.ClearFormats()
$ws.Range('E30').Value = '  -0.26%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.194'
$c.ClearFormats()
$ws.Range('E31').Value = '  -0.79%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.05254'
$c.ClearFormats()
$ws.Range('E32').Value = '  -0.69%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.847'
$c.ClearFormats()
$ws.Range('E33').Value = '  +0.62%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.7369'
$c.ClearFormats()
$ws.Range('E34').Value = '  -0.03%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.137'
$c.ClearFormats()
$ws.Range('E35').Value = '  -1.34%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.679'
$c.ClearFormats()
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('D37').Value = '1.246.71'
$ws.Range('E37').Value = '  -3.67%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.765'
$c.ClearFormats()
$ws.Range('E38').Value = '  +0.77%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01785'
$c.ClearFormats()
$ws.Range('E39').Value = '  +0.14%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '6.368'
$c.ClearFormats()
$ws.Range('E40').Value = '  +1.26%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.9000'
$c.ClearFormats()
$ws.Range('E41').Value = '  +0.43%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.ClearFormats()
$ws.Range('E42').Value = '  +0.06%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '102.02'
$c.ClearFormats()
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('D44').Value = '1.986.10'
$ws.Range('E44').Value = '  +0.57%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.00000000123'
$c.ClearFormats()
$ws.Range('E45').Value = '  +2.75%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '64.28'
$c.ClearFormats()
$ws.Range('E46').Value = '  +0.24%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.5121'
$c.ClearFormats()
$ws.Range('E47').Value = '  -0.29%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.4041'
$c.ClearFormats()
$ws.Range('E48').Value = '  +1.54%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '8.871'
$c.ClearFormats()
$ws.Range('E49').Value = '  +1.73%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.05775'
$c.ClearFormats()
$ws.Range('E50').Value = '  -0.47%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '6.719'
$c.ClearFormats()
$ws.Range('E51').Value = '  +0.35%  '
